$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C92").Value = "Average_Drawdown_test"
$ws.Range("B92").Value = "Test average drawdown"
$ws.Range("A92").Value = "Average Drawdown"

$ws.Range("A92").Select()
$excel.ActiveWindow.ScrollRow = 71
$excel.ActiveWindow.ScrollColumn = 1
